$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell text content per the commit:
#  B1: "windowUIType" -> "name"
#  M1: drop the "dialogueUIType, " prefix from the long comment text
$ws.Range("B1").Value = "name"
$ws.Range("M1").Value = "dialogDisplayType, buttonDisplayType are enum." + [char]10 + "widthRatio/heightRatio are from 0~1." + [char]10 + "If buttonText are not written, it will use the default text." + [char]10 + "title_EN is the title on the dialogBox." + [char]10 + "description_EN is the desctiption to the player."

# Move the active selection to M1 (matches the saved sheetView selection)
$ws.Range("M1").Select()
